# Auto-generated edit script applying the cryptos.xlsx cell updates
# described by the commit "Updated cryptos list on Mon Nov 13 21:43:05 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. NumberFormat is forced to Text ("@") before
# assignment so that numeric-looking strings (e.g. "241.98") are preserved as
# literal text, matching the original inlineStr cell type, instead of being
# auto-converted into Excel numbers (which would also lose exact formatting).
$cellValues = [ordered]@{
    'D2' = '36.524.05'
    'E2' = '  -1.81%  '
    'D3' = '2.062.71'
    'E3' = '  -0.04%  '
    'D4' = '0.999'
    'E4' = '  -0.14%  '
    'D5' = '241.98'
    'E5' = '  -2.76%  '
    'E6' = '  -1.21%  '
    'E7' = '  +0.03%  '
    'D8' = '52.70'
    'E8' = '  -7.89%  '
    'D9' = '58.76'
    'E9' = '  -2.62%  '
    'E10' = '  -7.65%  '
    'E11' = '  -5.28%  '
    'E12' = '  -0.17%  '
    'D13' = '0.892'
    'E13' = '  -2.68%  '
    'D14' = '14.64'
    'E14' = '  -10.02%  '
    'D15' = '2.360.78'
    'E15' = '  -0.01%  '
    'D16' = '5.40'
    'E16' = '  -6.76%  '
    'D17' = '2.094.37'
    'E17' = '  +1.52%  '
    'D18' = '36.412.07'
    'E18' = '  -2.15%  '
    'D19' = '16.43'
    'E19' = '  -12.08%  '
    'D20' = '71.38'
    'E20' = '  -4.77%  '
    'E21' = '  -5.30%  '
    'D22' = '5.26'
    'E22' = '  -4.59%  '
    'D23' = '236.16'
    'E23' = '  -0.80%  '
    'E24' = '  -0.01%  '
    'E25' = '  -5.18%  '
    'D26' = '9.41'
    'E26' = '  -3.07%  '
    'E27' = '  -3.06%  '
    'D28' = '163.98'
    'E28' = '  -3.56%  '
    'D29' = '20.27'
    'E29' = '  -0.21%  '
    'D30' = '0.122'
    'E30' = '  -2.55%  '
    'D31' = '5.07'
    'E31' = '  -2.15%  '
    'D32' = '1.13'
    'E32' = '  -2.84%  '
    'D33' = '4.58'
    'E33' = '  -1.71%  '
    'D34' = '0.0591'
    'E34' = '  -5.64%  '
    'E35' = '  +2.63%  '
    'E36' = '  -0.03%  '
    'E37' = '  +3.39%  '
    'D38' = '0.0808'
    'E38' = '  -9.38%  '
    'E39' = '  -8.21%  '
    'B40' = 'HuobiToken'
    'C40' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D40' = '2.93'
    'E40' = '  -6.44%  '
    'B41' = 'THORChain'
    'C41' = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
    'D41' = '4.81'
    'E41' = '  -9.06%  '
    'E42' = '  -2.91%  '
    'B43' = 'VeChain'
    'C43' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D43' = '0.0215'
    'E43' = '  -4.07%  '
    'D44' = '0.0939'
    'E44' = '  -7.85%  '
    'D45' = '93.67'
    'E45' = '  -3.22%  '
    'B46' = 'Maker'
    'C46' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D46' = '1.380.33'
    'E46' = '  +8.10%  '
    'B47' = 'FraxShare'
    'C47' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D47' = '7.47'
    'E47' = '  +8.98%  '
    'D48' = '15.25'
    'E48' = '  -13.85%  '
    'E49' = '  -4.37%  '
    'D50' = '2.85'
    'E50' = '  -0.32%  '
    'D51' = '2.249.31'
    'E51' = '  +0.02%  '
}

foreach ($addr in $cellValues.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cellValues[$addr]
}
